$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps for the
# d195734e-... report row (row 3) on both the zh-cn and de-de sheets,
# reflecting the regenerated handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 12:32:36"
$wsZhCn.Range("H3").Value = "2016-03-17 12:33:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 12:32:40"
$wsDeDe.Range("H3").Value = "2016-03-17 12:33:09"
